# "update column regarding changes"
# Insert a new "PMIDDLENAME" column right before the existing "PLASTNAME"
# column (column G) on Sheet1, shifting PLASTNAME and all the columns after
# it one slot to the right (H->I, I->J, J->K, K->L, L->M).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Put the selection back on A1 before we start (closest achievable match
# to the saved file no longer pointing the cursor at the old H8 cell).
$null = $ws.Activate()
$null = $ws.Range("A1").Select()

# Insert a new blank column at G - PLASTNAME (formerly G) and everything to
# its right shift one column over (G->H, H->I, I->J, J->K, K->L, L->M).
$ws.Columns("G").Insert()

# Header text for the newly-inserted column.
$ws.Range("G1").Value = "PMIDDLENAME"

# Give the new header cell the same look (bold header style) as the rest of
# row 1 by copying the format from the neighboring header cell.
$null = $ws.Range("F1").Copy()
$null = $ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Match the column width used for the new column.
$ws.Columns("G").ColumnWidth = 13.3

# Leave the selection cleared/back on A1.
$null = $ws.Range("A1").Select()
